# Generate Report for Handoff
# Updates the "Latest Handoff/HO Xliff Generate Date" timestamps and marks
# the Priority column as "ht" (handoff type) for the rows whose handoff
# report was just (re)generated.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn = $wb.Worksheets.Item("zh-cn")
$dede = $wb.Worksheets.Item("de-de")

$rows = @(7, 8, 10, 12, 13, 14)

foreach ($r in $rows) {
    # Overview sheet: column G = "Latest HO Xliff Generate Date"
    $overview.Range("G$r").Value = "2016-09-05 14:26:37"

    # zh-cn sheet: column E = "Priority", column H = "Latest Handoff Datetime"
    $zhcn.Range("E$r").Value = "ht"
    $zhcn.Range("H$r").Value = "2016-09-05 14:26:32"

    # de-de sheet: column E = "Priority", column H = "Latest Handoff Datetime"
    $dede.Range("E$r").Value = "ht"
    $dede.Range("H$r").Value = "2016-09-05 14:26:37"
}
